$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "SamplesTab" query (row 3, column B) so the Tumor column
# uses the sample's actual tumor status instead of the collected list.
$newQuery = "MATCH (s:study)<--(p:participant)<--(samp:sample)`n" +
    "WHERE s.study_name in [`"Detection of Colorectal Cancer Susceptibility Loci Using Genome-Wide Sequencing`"]`n" +
    "WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`n" +
    "RETURN  `n" +
    " coalesce(samp.sample_id, '') as ``Sample ID``,`n" +
    " coalesce(p.participant_id,'') as ``Participant ID``,`n" +
    " coalesce(s.study_name, '') as ``Study Name``,`n" +
    " coalesce(s.phs_accession,'') as ``Accession``,`n" +
    "coalesce(samp.sample_tumor_status,'') as ``Tumor``,`n" +
    "coalesce(samp.sample_type,'') as ``Analyte Type```n" +
    "ORDER By samp.sample_id LIMIT 100"

$ws.Range("B3").Value = $newQuery

# Move the active selection from C12 to C11.
$ws.Range("C11").Select()
